# Update monthly economic variables (variables_definitions.xlsx / Tabelle1)
#
# Adds six new variable rows to the lookup table:
#   1. RetTurn   - Retail turnover excluding cars          (Activity group)
#   2. EPI       - Export price index                      (Prices group)
#   3. IPI       - Import price index                      (Prices group)
#   4. Empl      - Employment                               (Labor market group)
#   5. GWMan     - Gross wages and salaries: manufacturing and mining (Labor market)
#   6. GWConstr  - Gross wages and salaries: construction   (Labor market group)
#
# The whole table (including previously-existing rows) is rewritten in place
# so that the final table is ordered by Group (Activity / Prices / Labor
# market) with the new rows inserted at the end of their respective group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Mnemonic"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Group"
$ws.Range("A2").Value = "ConstrProd"
$ws.Range("B2").Value = "Production in main construction industry"
$ws.Range("C2").Value = "Activity"
$ws.Range("A3").Value = "IP"
$ws.Range("B3").Value = "Industrial production index"
$ws.Range("C3").Value = "Activity"
$ws.Range("A4").Value = "ConstrNO"
$ws.Range("B4").Value = "New orders for main construction industry"
$ws.Range("C4").Value = "Activity"
$ws.Range("A5").Value = "INO"
$ws.Range("B5").Value = "New orders for industry"
$ws.Range("C5").Value = "Activity"
$ws.Range("A6").Value = "ConstrTurn"
$ws.Range("B6").Value = "Main construction industry turnover"
$ws.Range("C6").Value = "Activity"
$ws.Range("A7").Value = "ITurn"
$ws.Range("B7").Value = "Industry turnover"
$ws.Range("C7").Value = "Activity"
$ws.Range("A8").Value = "RetTurn"
$ws.Range("B8").Value = "Retail turnover excluding cars"
$ws.Range("C8").Value = "Activity"
$ws.Range("A9").Value = "CPI"
$ws.Range("B9").Value = "Consumer price index"
$ws.Range("C9").Value = "Prices"
$ws.Range("A10").Value = "CPIEN"
$ws.Range("B10").Value = "Consumer price index, excluding energy"
$ws.Range("C10").Value = "Prices"
$ws.Range("A11").Value = "PPI"
$ws.Range("B11").Value = "Producer price index"
$ws.Range("C11").Value = "Prices"
$ws.Range("A12").Value = "PPIEN"
$ws.Range("B12").Value = "Producer price index, excluding energy"
$ws.Range("C12").Value = "Prices"
$ws.Range("A13").Value = "EPI"
$ws.Range("A14").Value = "IPI"
$ws.Range("B13").Value = "Export price index"
$ws.Range("B14").Value = "Import price index"
$ws.Range("C13").Value = "Prices"
$ws.Range("C14").Value = "Prices"
$ws.Range("A15").Value = "HW"
$ws.Range("B15").Value = "Hours worked: manufacturing"
$ws.Range("C15").Value = "Labor market"
$ws.Range("A16").Value = "ConstrHW"
$ws.Range("B16").Value = "Hours worked: construction"
$ws.Range("C16").Value = "Labor market"
$ws.Range("A17").Value = "Empl"
$ws.Range("B17").Value = "Employment"
$ws.Range("C17").Value = "Labor market"
$ws.Range("A18").Value = "GWMan"
$ws.Range("B18").Value = "Gross wages and salaries: manufacturing and mining"
$ws.Range("C18").Value = "Labor market"
$ws.Range("A19").Value = "GWConstr"
$ws.Range("B19").Value = "Gross wages and salaries: construction"
$ws.Range("C19").Value = "Labor market"
$ws.Range("C19").Select()
